# Revert "Merge branch 'master' of https://github.com/pdisbeschl/Uni_Timetables"
# This reverts the lecturer-schedule edits that were merged in, restoring the
# "Lecturers" sheet to its prior content (Kaestner/Paul rows, no R.Möckel row)
# and restoring the previously-active sheet/selection state.

$wb = $excel.ActiveWorkbook
$wsL = $wb.Worksheets.Item("Lecturers")

# --- Lecturers sheet content ---

# Rows 2-4 used to be "P. Bonizzi" (style of row1/A1 plain) - revert them back
# to "Kaestner" using the same formatting already used by row 5 (A5).
$wsL.Range("A5").Copy()
$wsL.Range("A2:A4").PasteSpecial(-4122)

$wsL.Range("A2").Value = "Kaestner"
$wsL.Range("A3").Value = "Kaestner"
$wsL.Range("A4").Value = "Kaestner"

# Restore the previous date values for each lecturer session.
$wsL.Range("B2").Value = 43941
$wsL.Range("B3").Value = 43942
$wsL.Range("B4").Value = 43943
$wsL.Range("B5").Value = 43944
$wsL.Range("B6").Value = 43945
$wsL.Range("B7").Value = 43946
$wsL.Range("B8").Value = 43947

# Restore the previous end-time value for row 4.
$wsL.Range("D4").Value = 0.54166666666666663

# The extra "R.Möckel" row (row 9) did not exist before - remove it.
$wsL.Rows.Item(9).Delete()

# --- Restore previous selections / active sheet ---

$wb.Worksheets.Item("Courses").Range("E23").Select()
$wsL.Range("C4").Select()
$wb.Worksheets.Item("Holidays").Select()
